$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.187699999999991
$ws.Range("B4").Value = 5.634700000000002
$ws.Range("D6").Value = -7.960399999999997
$ws.Range("B7").Value = 5.909099999999992
$ws.Range("D7").Value = -7.092599999999994
$ws.Range("B8").Value = 6.142799999999997
$ws.Range("D8").Value = -7.6939
$ws.Range("A11").Value = -21.91580000000002
$ws.Range("E11").Value = 13.41710000000001
$ws.Range("A12").Value = -22.72080000000002
$ws.Range("B12").Value = 5.7308
$ws.Range("B14").Value = 8.978900000000003
$ws.Range("E14").Value = 13.03390000000001
$ws.Range("A15").Value = -21.35760000000003
$ws.Range("D19").Value = -8.118299999999991
$ws.Range("E19").Value = 13.5879
$ws.Range("D21").Value = -7.748000000000003
$ws.Range("E21").Value = 13.19270000000002
$ws.Range("B22").Value = 5.8798
$ws.Range("D24").Value = -8.169199999999996
$ws.Range("D25").Value = -7.866199999999992
